# Reorder the "readme" sheet's log table columns from
#   index, Date, Author, sheet_name, JobNo
# to
#   index, Author, JobNo, Date, sheet_name
# and refresh the recorded "Date of Analysis" timestamp on the
# "Project Information" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("readme")

# --- Header row (row 1) --------------------------------------------------
# Plain words, so assigning .Value directly is safe (no numeric
# reinterpretation) and keeps the Table1 ListObject column names in sync.
$ws.Range("B1").Value = "Author"
$ws.Range("C1").Value = "JobNo"
$ws.Range("D1").Value = "Date"
$ws.Range("E1").Value = "sheet_name"

# --- Data rows (rows 2-12) ------------------------------------------------
# The "Date" column holds digit-only text ("20220224"); assigning it with
# .Value would make Excel reinterpret it as a number. Instead, rotate the
# four columns with Copy (through a scratch area) so each cell's original
# text typing and style travel with it unchanged.
#
# New layout per column, in terms of the old column holding that data:
#   B(Author) <- old C, C(JobNo) <- old E, D(Date) <- old B, E(sheet_name) <- old D
$ws.Range("B2:B12").Copy($ws.Range("G2"))   # stash old Date
$ws.Range("C2:C12").Copy($ws.Range("H2"))   # stash old Author
$ws.Range("D2:D12").Copy($ws.Range("I2"))   # stash old sheet_name
$ws.Range("E2:E12").Copy($ws.Range("J2"))   # stash old JobNo

$ws.Range("H2:H12").Copy($ws.Range("B2"))   # B = Author
$ws.Range("J2:J12").Copy($ws.Range("C2"))   # C = JobNo
$ws.Range("G2:G12").Copy($ws.Range("D2"))   # D = Date
$ws.Range("I2:I12").Copy($ws.Range("E2"))   # E = sheet_name

$ws.Range("G2:J12").Clear()

# --- Refresh analysis timestamp ------------------------------------------
$ws2 = $wb.Worksheets.Item("Project Information")
$ws2.Range("B11").Value = "2022-02-24 13:15:49.867676"
